{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\n\n// 1) Table cell text change: \"MIERCOLES\" -> \"LUNES\"\nconst dayResults = body.search(\"MIERCOLES\", { matchCase: true, matchWholeWord: true });\ndayResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < dayResults.items.length; i++) {\n  dayResults.items[i].insertText(\"LUNES\", Word.InsertLocation.replace);\n}\n\n// 2) Merge the two runs around the removed \"_GoBack\" bookmark back into a\n//    single run (\"...ambas partes\" + \".\" -> \"...ambas partes.\") and drop\n//    the now-unused bookmark.\nconst tailResults = body.search(\"ambas partes.\", { matchCase: true });\nawait context.sync();\n\nif (tailResults.items.length > 0) {\n  // Re-inserting the same text over the matched range merges the runs that\n  // previously straddled the bookmark into one run, preserving formatting.\n  tailResults.items[0].insertText(\"ambas partes.\", Word.InsertLocation.replace);\n}\n\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# 1) Table cell text change: \"MIERCOLES\" -> \"LUNES\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"MIERCOLES\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"LUNES\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)\n\n# 2) Remove the now-unused \"_GoBack\" bookmark that split the closing\n#    sentence into two runs.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 3) Re-apply the trailing text over the old bookmark boundary so the two\n#    runs (\"...ambas partes\" + \".\") merge back into a single run.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"ambas partes.\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"ambas partes.\"\n$find2.MatchCase = $true\n$find2.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)\n"}
